$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(1.551991501684256,1.37776271018447,1.271467160817167,1.228315569744439,1.221160035919183,1.270884545222145,1.49177182103989,1.930668248630298,2.25711722778982,2.406605810174312,2.463363635692303,2.451133031878612,2.411272262603973,2.386876134457168,2.24736830803397,2.162042799645235,2.113058070332499,2.096488284122017,2.171116240855895,2.422976207460863,2.588457975055405,2.500054300966269,2.167013920525392,1.811270555958458)
$colC = @(0.6356254217462265,0.5565799236435964,0.50828578767198,0.4886617989130286,0.4854065297942043,0.5080209092781729,0.6083181766624648,0.8071035520439409,0.9547238745029176,1.022286925098285,1.047935406917645,1.042408629878139,1.024395733631081,1.013370773093243,0.9503171321304649,0.9117441722768831,0.8895963011008803,0.8821038392931655,0.9158463333425289,1.029684786259168,1.104459148958256,1.064514764078524,0.9139916587402581,0.753072052581274)
$colE = @(0.4189517054670944,0.3649226721706498,0.3319174917557319,0.3185061684469304,0.316281444221346,0.3317364709654385,0.4002850610186499,0.5362332550958797,0.6373123177370275,0.68361605766016,0.7012013049098442,0.6974116622608335,0.6850617567150294,0.6775038687183752,0.6342931427580538,0.6078706130607259,0.5927031290883349,0.5875727015139347,0.6106801922959306,0.6886878007014161,0.7399699545040619,0.7125707092548481,0.609409909673829,0.4992650714670361)
$colF = @(0.4443680307746121,0.3878228170618172,0.3531389305169483,0.3390132514313251,0.336668177824194,0.3529483938368969,0.4248636149813905,0.5661985755042025,0.6702781546542269,0.7176906081379002,0.7356546913071611,0.7317853510981394,0.7191683204515869,0.7114413442032514,0.6671810134426437,0.6400460337215605,0.6244449056556647,0.6191636801734006,0.6429339538360921,0.7228739723492197,0.7751780083420101,0.7472568307916134,0.6416283278902171,0.5279251897347308)
$colG = @(0.002405939069815402,0.002412601349006005,0.002416892817716384,0.002418692344663421,0.002418994224562265,0.002416916881149708,0.002408194683528328,0.002392673408746092,0.002382220330240365,0.002377668169104021,0.002375973328281821,0.002376337057931198,0.002377528154304424,0.002378261500348605,0.002382521890145395,0.002385187332382827,0.002386739545024962,0.002387268388506062,0.002384901614317079,0.002377177516107443,0.002372298089386293,0.002374886968895823,0.002385030725570155,0.002396704366730327)
$colI = @(0.7585760744539272,0.7375074338506096,0.7252513311404556,0.7204251029660398,0.719633791474493,0.7251855654249155,0.7511687474827511,0.8076488590455853,0.8527023129189644,0.8740126251813223,0.8822029268833944,0.8804335838799204,0.8746840110158161,0.871178030222552,0.8513263208980533,0.8393590862122977,0.8325523840099009,0.8302607966581519,0.8406250739483454,0.8763695015060478,0.9004347920796079,0.8875251417975534,0.8400524924544897,0.791758563825681)
$colN = @(1.02950034710998,1.044842132735837,1.054810370801782,1.05900933531106,1.059714803632211,1.054866446874875,1.034675754878535,0.9994797420445778,0.9763672942018147,0.966463499868297,0.9628021839106822,0.9635867329948624,0.9661604892784084,0.9677486233058801,0.9770269186641514,0.9828760441896875,0.9862976599282405,0.9874659717485272,0.9822474496504654,0.9654020878037457,0.9549122607224021,0.9604629013972996,0.9825314539549126,1.008523746835372)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, "B").Value = $colB[$i]
    $ws.Cells.Item($row, "C").Value = $colC[$i]
    $ws.Cells.Item($row, "E").Value = $colE[$i]
    $ws.Cells.Item($row, "F").Value = $colF[$i]
    $ws.Cells.Item($row, "G").Value = $colG[$i]
    $ws.Cells.Item($row, "I").Value = $colI[$i]
    $ws.Cells.Item($row, "N").Value = $colN[$i]
}